$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-13 04:42:38"
$wsZhCn.Range("H3").Value = "2016-03-13 04:42:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-13 04:42:42"
$wsDeDe.Range("H3").Value = "2016-03-13 04:43:04"
